$wb = $excel.ActiveWorkbook

# Hunk 0: @@ -874,25 +874,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 245.78572
$ws.Range("I5").Value = 82.44444
$ws.Range("J5").Value = 539.8
$ws.Range("K5").Value = 82.44444
$ws.Range("L5").Value = 539.8
$ws.Range("M5").Value = 32.55556
$ws.Range("N5").Value = -769.8

# Hunk 1: @@ -7532,25 +7532,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2585.9136
$ws.Range("I138").Value = 1340.8125
$ws.Range("J138").Value = 3399.0408
$ws.Range("K138").Value = 4022.4375
$ws.Range("L138").Value = 10197.1224
$ws.Range("M138").Value = 1117.5625
$ws.Range("N138").Value = -20477.1224

# Hunk 2: @@ -7828,22 +7828,22 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 155241.77
$ws.Range("I2").Value = 223515.89
$ws.Range("K2").Value = 223515.89
$ws.Range("M2").Value = -223402.89

# Hunk 3: @@ -9858,25 +9858,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 9339.75
$ws.Range("I43").Value = 6121.3335
$ws.Range("J43").Value = 11270.8
$ws.Range("K43").Value = 6121.3335
$ws.Range("L43").Value = 11270.8
$ws.Range("M43").Value = -5808.3335
$ws.Range("N43").Value = -11896.8

# Hunk 4: @@ -11383,25 +11383,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1777.1034
$ws.Range("I74").Value = 1490.9166
$ws.Range("J74").Value = 1979.1177
$ws.Range("K74").Value = 1490.9166
$ws.Range("L74").Value = 1979.1177
$ws.Range("M74").Value = -616.9166
$ws.Range("N74").Value = -3727.1177

# Hunk 5: @@ -11530,25 +11530,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1777.1034
$ws.Range("I77").Value = 1490.9166
$ws.Range("J77").Value = 1979.1177
$ws.Range("K77").Value = 7454.583000000001
$ws.Range("L77").Value = 9895.5885
$ws.Range("M77").Value = -3086.583000000001
$ws.Range("N77").Value = -18631.5885

# Hunk 6: @@ -13432,22 +13432,22 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 155241.77
$ws.Range("I116").Value = 223515.89
$ws.Range("K116").Value = 223515.89
$ws.Range("M116").Value = -221221.89

# Hunk 7: @@ -13723,25 +13723,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1512730.9
$ws.Range("I122").Value = 2141785.2
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6425355.600000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6422905.600000001
$ws.Range("N122").Value = -13900

# Hunk 8: @@ -14852,22 +14852,22 @@ (sheet BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 155241.77
$ws.Range("I3").Value = 223515.89
$ws.Range("K3").Value = 223515.89
$ws.Range("M3").Value = -223401.89

# Hunk 9: @@ -22398,25 +22398,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1851.75
$ws.Range("I16").Value = 1343.2222
$ws.Range("J16").Value = 2505.5715
$ws.Range("K16").Value = 1343.2222
$ws.Range("L16").Value = 2505.5715
$ws.Range("M16").Value = -1056.2222
$ws.Range("N16").Value = -3079.5715

# Hunk 10: @@ -26291,19 +26291,22 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 38500
$ws.Range("J96").Value = 38500
$ws.Range("L96").Value = 38500
$ws.Range("N96").Value = -43992

# Hunk 11: @@ -27118,25 +27121,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1851.75
$ws.Range("I113").Value = 1343.2222
$ws.Range("J113").Value = 2505.5715
$ws.Range("K113").Value = 1343.2222
$ws.Range("L113").Value = 2505.5715
$ws.Range("M113").Value = 826.7778000000001
$ws.Range("N113").Value = -6845.5715

# Hunk 12: @@ -28034,25 +28037,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2358.25
$ws.Range("I132").Value = 1699.5
$ws.Range("J132").Value = 3393.4285
$ws.Range("K132").Value = 5098.5
$ws.Range("L132").Value = 10180.2855
$ws.Range("M132").Value = -2568.5
$ws.Range("N132").Value = -15240.2855

# Hunk 13: @@ -33429,25 +33432,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 9091528
$ws.Range("I97").Value = 12500639
$ws.Range("J97").Value = 566.6667
$ws.Range("K97").Value = 37501917
$ws.Range("L97").Value = 1700.0001
$ws.Range("M97").Value = -37501421
$ws.Range("N97").Value = -2692.0001

# Hunk 14: @@ -34234,25 +34237,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3333841.5
$ws.Range("I113").Value = 6250435
$ws.Range("J113").Value = 1000566.8
$ws.Range("K113").Value = 18751305
$ws.Range("L113").Value = 3001700.4
$ws.Range("M113").Value = -18749135
$ws.Range("N113").Value = -3006040.4

# Hunk 15: @@ -39646,22 +39649,22 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16429.285
$ws.Range("I80").Value = 22001
$ws.Range("K80").Value = 22001
$ws.Range("M80").Value = -21003

# Hunk 16: @@ -39793,22 +39796,22 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 16429.285
$ws.Range("I83").Value = 22001
$ws.Range("K83").Value = 110005
$ws.Range("M83").Value = -105013

# Hunk 17: @@ -41668,25 +41671,25 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 262963760
$ws.Range("I122").Value = 266204700
$ws.Range("J122").Value = 250000000
$ws.Range("K122").Value = 798614100
$ws.Range("L122").Value = 750000000
$ws.Range("M122").Value = -798611650
$ws.Range("N122").Value = -750004900

# Hunk 18: @@ -42978,25 +42981,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3378.889
$ws.Range("I7").Value = 1850
$ws.Range("J7").Value = 6436.6665
$ws.Range("K7").Value = 1850
$ws.Range("L7").Value = 6436.6665
$ws.Range("M7").Value = -1738
$ws.Range("N7").Value = -6660.6665

# Hunk 19: @@ -43707,25 +43710,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2333.2856
$ws.Range("I22").Value = 1600
$ws.Range("J22").Value = 2455.5
$ws.Range("K22").Value = 1600
$ws.Range("L22").Value = 2455.5
$ws.Range("M22").Value = -1305
$ws.Range("N22").Value = -3045.5

# Hunk 20: @@ -43955,25 +43958,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2333.2856
$ws.Range("I27").Value = 1600
$ws.Range("J27").Value = 2455.5
$ws.Range("K27").Value = 1600
$ws.Range("L27").Value = 2455.5
$ws.Range("M27").Value = -1493
$ws.Range("N27").Value = -2669.5

# Hunk 21: @@ -44592,23 +44595,26 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 500001100
$ws.Range("J40").Value = 2180
$ws.Range("L40").Value = 2180
$ws.Range("N40").Value = -2452

# Hunk 22: @@ -45321,22 +45327,22 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 154
$ws.Range("I55").Value = 250
$ws.Range("K55").Value = 250
$ws.Range("M55").Value = -77

# Hunk 23: @@ -47383,22 +47389,22 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 36500
$ws.Range("J98").Value = 36500
$ws.Range("L98").Value = 36500
$ws.Range("N98").Value = -42490

# Hunk 24: @@ -48523,25 +48529,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 16296614
$ws.Range("I122").Value = 23826924
$ws.Range("J122").Value = 5001150
$ws.Range("K122").Value = 71480772
$ws.Range("L122").Value = 15003450
$ws.Range("M122").Value = -71478322
$ws.Range("N122").Value = -15008350

# Hunk 25: @@ -48722,25 +48728,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3378.889
$ws.Range("I126").Value = 1850
$ws.Range("J126").Value = 6436.6665
$ws.Range("K126").Value = 5550
$ws.Range("L126").Value = 19309.9995
$ws.Range("M126").Value = -3080
$ws.Range("N126").Value = -24249.9995

# Hunk 26: @@ -51487,22 +51493,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28843.5
$ws.Range("J41").Value = 28843.5
$ws.Range("L41").Value = 28843.5
$ws.Range("N41").Value = -29623.5

# Hunk 27: @@ -51680,22 +51686,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 6875.3335
$ws.Range("J45").Value = 6875.3335
$ws.Range("L45").Value = 6875.3335
$ws.Range("N45").Value = -7857.3335

# Hunk 28: @@ -53092,22 +53098,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 36346.668
$ws.Range("J74").Value = 36346.668
$ws.Range("L74").Value = 36346.668
$ws.Range("N74").Value = -38218.668

# Hunk 29: @@ -53236,22 +53242,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 36346.668
$ws.Range("J77").Value = 36346.668
$ws.Range("L77").Value = 109040.004
$ws.Range("N77").Value = -118400.004

# Hunk 30: @@ -54544,22 +54550,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 40074
$ws.Range("J104").Value = 40074
$ws.Range("L104").Value = 40074
$ws.Range("N104").Value = -47062

# Hunk 31: @@ -55402,25 +55408,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 795
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 2385
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -7285
